# Add payment 09876543 (Cash) 2025-08-18T18:03:23
#
# Row 52 was the most-recently appended payment; its phone number had been
# stored as text ("09876543", with the leading zero). Recording the new
# payment normalizes that earlier cell to a plain number, and the new
# payment itself is appended as row 53 with the phone number kept as text
# (so the leading zero survives).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Normalize the previous row's phone number to a number ---
$ws.Range("A52").Value = 9876543

# --- Append the new payment row (53) ---

# phone (keep the leading zero -> must stay text)
$ws.Range("A53").NumberFormat = "@"
$ws.Range("A53").Value = "09876543"
$ws.Range("A53").ClearFormats()

# method
$ws.Range("C53").Value = "Cash"

# timestamp (stored as text, not an Excel date)
$ws.Range("D53").NumberFormat = "@"
$ws.Range("D53").Value = "2025-08-18T18:03:23"
$ws.Range("D53").ClearFormats()

# original_amount / final_amount / birthday_discount / points_redeemed / reward_discount
$ws.Range("E53").Value = 120
$ws.Range("G53").Value = 120
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("J53").Value = 0
